$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 87 ---
# Only the date/time value in column A changes; everything else stays the same.
$ws.Range("A87").Value = 45468.2916666667

# --- Append new row 88 ---
$ws.Range("A88").Value = 45469.6002777778
$ws.Range("B88").Value = 3000
$ws.Range("C88").Value = 3.07999992370605
$ws.Range("D88").Value = 3.01999998092651
$ws.Range("E88").Value = 3.01999998092651
$ws.Range("F88").Value = 3.07999992370605
$ws.Range("H88").Value = "ESPE.MI"

# G88 must be stored as a shared-string "3.07999992370605" (matching G-column
# convention for adj_close) rather than a numeric value. Direct assignment of a
# numeric-looking string gets auto-converted to a number by Excel, so build the
# text in a scratch cell via a formula (forces text type), then copy/paste the
# *value* into G88 and clean up the scratch cell.
$ws.Range("Z1").Formula = "=""3.07999992370605"""
$ws.Range("Z1").Copy()
$ws.Range("G88").PasteSpecial(-4163) # xlPasteValues
$ws.Range("Z1").ClearContents()

# Give A88 the same date/time number format style as the rest of column A.
$ws.Range("A87").Copy()
$ws.Range("A88").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false
